## revised designation rec paragraph
## - Adds a "Designation Recommendation Text" paragraph (with its own new
##   paragraph style) right after the "Designation Recommendation" heading
##   paragraph.
## - Adds w:spacing w:before="360" to the "Designation Recommendation" style.
## - Bumps the cached "July 2025" DATE field text to "August 2025" in both
##   the default and first-page footers.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Footer date field text: "July 2025" -> "August 2025"
# ---------------------------------------------------------------------
$sec = $d.Sections.Item(1)
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $footer = $sec.Footers.Item($i)
    if ($footer.Exists) {
        $footer.Range.Find.Execute("July 2025", $true, $false, $false, $false, $false, $true, 1, $false, "August 2025", 2) | Out-Null
    }
}

# ---------------------------------------------------------------------
# 2. "Designation Recommendation" style: add spacing-before of 18pt (360 twips)
# ---------------------------------------------------------------------
$designationStyle = $d.Styles.Item("Designation Recommendation")
$designationStyle.ParagraphFormat.SpaceBefore = 18

# ---------------------------------------------------------------------
# 3. New paragraph style "Designation Recommendation Text", based on
#    Body Text, with left/right indents of 720 twips (36pt) each.
# ---------------------------------------------------------------------
$textStyle = $d.Styles.Add("Designation Recommendation Text", 1)
$textStyle.BaseStyle = "BodyText"
$textStyle.QuickStyle = $true
$textStyle.ParagraphFormat.LeftIndent = 36
$textStyle.ParagraphFormat.RightIndent = 36

# ---------------------------------------------------------------------
# 4. Insert the new "Designation Recommendation Text" paragraph right
#    after the "Designation Recommendation" paragraph, with its runs
#    (one of them carrying the VerbatimChar character style).
# ---------------------------------------------------------------------
$designationPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Designation Recommendation") {
        $designationPara = $p
    }
}
if ($designationPara -eq $null) {
    throw "Could not find the 'Designation Recommendation' paragraph."
}

$designationPara.Range.InsertParagraphAfter()
$newPara = $designationPara.Next()
$newParaRange = $newPara.Range

$newParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="DesignationRecommendationText"/></w:pPr><w:r><w:t xml:space="preserve">Designation Recommendation Text: </w:t></w:r><w:r><w:t xml:space="preserve">Note that the </w:t></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t>echo = FALSE</w:t></w:r><w:r><w:t xml:space="preserve"> parameter was added to the code chunk to prevent printing of the R code that generated the plot.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId2" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/styles" Target="styles.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/styles.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.styles+xml"><pkg:xmlData><w:styles xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:style w:type="character" w:customStyle="1" w:styleId="VerbatimChar"><w:name w:val="Verbatim Char"/></w:style><w:style w:type="paragraph" w:customStyle="1" w:styleId="DesignationRecommendationText"><w:name w:val="Designation Recommendation Text"/></w:style></w:styles></pkg:xmlData></pkg:part></pkg:package>
'@

$newParaRange.InsertXML($newParaXml)
